$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.805.97'
$ws.Range('E2').Value = '  -2.01%  '
$ws.Range('D3').Value = '3.051.85'
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'556.73"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').Value = "'142.01"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.049.13'
$ws.Range('E8').Value = '  -1.78%  '
$ws.Range('D9').Value = "'0.519"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.61%  '
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('D11').Value = "'6.15"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.47%  '
$ws.Range('D12').Value = "'0.479"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.55%  '
$ws.Range('D13').Value = "'0.0000231"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('D14').Value = "'35.17"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').Value = '3.552.30'
$ws.Range('E15').Value = '  -1.50%  '
$ws.Range('D16').Value = '63.833.61'
$ws.Range('E16').Value = '  -1.90%  '
$ws.Range('D17').Value = '3.047.22'
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').Value = "'495.06"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('D21').Value = "'14.18"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.78%  '
$ws.Range('D22').Value = "'0.684"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').Value = "'14.61"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.47%  '
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').Value = "'82.69"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = "'8.12"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('E29').Value = '  -1.47%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = "'26.46"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').Value = "'2.46"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('D34').Value = "'5.69"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('D35').Value = "'6.20"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = "'55.38"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').Value = "'0.0409"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').Value = "'440.28"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.24%  '
$ws.Range('D39').Value = "'0.0815"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.14%  '
$ws.Range('D40').Value = '3.040.60'
$ws.Range('E40').Value = '  +1.27%  '
$ws.Range('D41').Value = "'2.77"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.90%  '
$ws.Range('D42').Value = "'8.33"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('D44').Value = "'0.274"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.56%  '
$ws.Range('D45').Value = "'27.67"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.38%  '
$ws.Range('D46').Value = "'2.23"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.14%  '
$ws.Range('D48').Value = "'0.114"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('D49').Value = "'117.89"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.22%  '
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('D51').Value = "'2.10"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.91%  '
